$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2039.2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2039.2
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6117.6
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -6453.6
$ws.Range("H33").Value = 132.5
$ws.Range("I33").Value = 124
$ws.Range("J33").Value = 175
$ws.Range("K33").Value = 124
$ws.Range("L33").Value = 175
$ws.Range("M33").Value = 105
$ws.Range("N33").Value = -633
$ws.Range("H43").Value = 3149.5
$ws.Range("I43").Value = 3066.3333
$ws.Range("K43").Value = 3066.3333
$ws.Range("M43").Value = -2997.3333
$ws.Range("H62").Value = 4999
$ws.Range("I62").Value = 4999
$ws.Range("K62").Value = 4999
$ws.Range("M62").Value = -4375
$ws.Range("H65").Value = 4999
$ws.Range("I65").Value = 4999
$ws.Range("K65").Value = 24995
$ws.Range("M65").Value = -21875
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H82").Value = 44
$ws.Range("I82").Value = 44
$ws.Range("K82").Value = 132
$ws.Range("M82").Value = 274
$ws.Range("H85").Value = 44
$ws.Range("I85").Value = 44
$ws.Range("K85").Value = 132
$ws.Range("M85").Value = 1272
$ws.Range("H130").Value = 95978
$ws.Range("J130").Value = 95978
$ws.Range("L130").Value = 95978
$ws.Range("N130").Value = -106018

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H32").Value = 1368.5
$ws.Range("I32").Value = 1532.8572
$ws.Range("K32").Value = 1532.8572
$ws.Range("M32").Value = -1245.8572
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H102").Value = 3271.2856
$ws.Range("I102").Value = 2780
$ws.Range("J102").Value = 4499.5
$ws.Range("K102").Value = 2780
$ws.Range("L102").Value = 4499.5
$ws.Range("M102").Value = -1158
$ws.Range("N102").Value = -7743.5
$ws.Range("H130").Value = 92776.336
$ws.Range("J130").Value = 92776.336
$ws.Range("L130").Value = 92776.336
$ws.Range("N130").Value = -102816.336
$ws.Range("H134").Value = 95000
$ws.Range("J134").Value = 95000
$ws.Range("L134").Value = 95000
$ws.Range("N134").Value = -105140
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H86").Value = 1073.75
$ws.Range("J86").Value = 1030
$ws.Range("L86").Value = 1030
$ws.Range("N86").Value = -3276
$ws.Range("H89").Value = 1073.75
$ws.Range("J89").Value = 1030
$ws.Range("L89").Value = 5150
$ws.Range("N89").Value = -16382
$ws.Range("H116").Value = 90000
$ws.Range("J116").Value = 90000
$ws.Range("L116").Value = 90000
$ws.Range("N116").Value = -99178
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H31").Value = 3618.5925
$ws.Range("J31").Value = 4714.9443
$ws.Range("L31").Value = 4714.9443
$ws.Range("N31").Value = -5304.9443
$ws.Range("H34").Value = 3618.5925
$ws.Range("J34").Value = 4714.9443
$ws.Range("L34").Value = 4714.9443
$ws.Range("N34").Value = -5118.9443
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H122").Value = 3965.6667
$ws.Range("J122").Value = 3965.6667
$ws.Range("L122").Value = 11897.0001
$ws.Range("N122").Value = -16797.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1183.3334
$ws.Range("J23").Value = 1237.5
$ws.Range("L23").Value = 3712.5
$ws.Range("N23").Value = -4182.5
$ws.Range("H34").Value = 1417.6666
$ws.Range("I34").Value = 750
$ws.Range("J34").Value = 1751.5
$ws.Range("K34").Value = 2250
$ws.Range("L34").Value = 5254.5
$ws.Range("M34").Value = -2166
$ws.Range("N34").Value = -5422.5
$ws.Range("H121").Value = 602.1429000000001
$ws.Range("I121").Value = 631.25
$ws.Range("J121").Value = 563.3333
$ws.Range("K121").Value = 1893.75
$ws.Range("L121").Value = 1689.9999
$ws.Range("M121").Value = -583.75
$ws.Range("N121").Value = -4309.9999
$ws.Range("H128").Value = 150000
$ws.Range("I128").Value = 150000
$ws.Range("K128").Value = 450000
$ws.Range("M128").Value = -445020

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 6800833.5
$ws.Range("I14").Value = 7315143
$ws.Range("J14").Value = 5000750
$ws.Range("K14").Value = 7315143
$ws.Range("L14").Value = 5000750
$ws.Range("M14").Value = -7314975
$ws.Range("N14").Value = -5001086
$ws.Range("H110").Value = 99959
$ws.Range("J110").Value = 99959
$ws.Range("L110").Value = 99959
$ws.Range("N110").Value = -108139
$ws.Range("H132").Value = 2550
$ws.Range("I132").Value = 2575
$ws.Range("K132").Value = 7725
$ws.Range("M132").Value = -5195
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140
$ws.Range("H138").Value = 84425
$ws.Range("J138").Value = 84425
$ws.Range("L138").Value = 84425
$ws.Range("N138").Value = -94705

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 19016.666
$ws.Range("I5").Value = 13500
$ws.Range("K5").Value = 13500
$ws.Range("M5").Value = -13387
$ws.Range("H24").Value = 15953
$ws.Range("I24").Value = 12006
$ws.Range("K24").Value = 12006
$ws.Range("M24").Value = -11663
$ws.Range("H29").Value = 19400
$ws.Range("I29").Value = 10000
$ws.Range("K29").Value = 10000
$ws.Range("M29").Value = -9705
$ws.Range("H40").Value = 4185.923
$ws.Range("J40").Value = 4300.1665
$ws.Range("L40").Value = 4300.1665
$ws.Range("N40").Value = -4572.1665
$ws.Range("H46").Value = 4834.3887
$ws.Range("I46").Value = 4000.5
$ws.Range("J46").Value = 4938.625
$ws.Range("K46").Value = 4000.5
$ws.Range("L46").Value = 4938.625
$ws.Range("M46").Value = -3812.5
$ws.Range("N46").Value = -5314.625
$ws.Range("H93").Value = 197.5
$ws.Range("I93").Value = 197.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 197.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 1050.5
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 4185.3335
$ws.Range("I100").Value = 4185.3335
$ws.Range("K100").Value = 4185.3335
$ws.Range("M100").Value = -3644.3335
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H134").Value = 95000
$ws.Range("J134").Value = 95000
$ws.Range("L134").Value = 95000
$ws.Range("N134").Value = -105140
$ws.Range("H136").Value = 4078.125
$ws.Range("I136").Value = 3803.5715
$ws.Range("K136").Value = 11410.7145
$ws.Range("M136").Value = -8860.7145
$ws.Range("H137").Value = 51500
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10000000
$ws.Range("J5").Value = 10000000
$ws.Range("L5").Value = 10000000
$ws.Range("N5").Value = -10000224
$ws.Range("H22").Value = 8006.5
$ws.Range("I22").Value = 6013
$ws.Range("K22").Value = 6013
$ws.Range("M22").Value = -5720
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

Write-Host "Applied all changes"